# Applies the price/volume refresh described in the commit diff.
# Prefixing each literal with an apostrophe forces text storage so values
# such as "1.002" or subscript-digit prices are not reinterpreted as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"  # text-prefix marker (literal apostrophe, kept as a string)

$ws.Range("D2").Value = $q + '25.658.43'
$ws.Range("E2").Value = $q + '  -1.77%  '
$ws.Range("D3").Value = $q + '1.615.29'
$ws.Range("E3").Value = $q + '  -1.90%  '
$ws.Range("D4").Value = $q + '1.002'
$ws.Range("E4").Value = $q + '  -0.12%  '
$ws.Range("D5").Value = $q + '214.14'
$ws.Range("E5").Value = $q + '  -1.26%  '
$ws.Range("D6").Value = $q + '0.5057'
$ws.Range("E6").Value = $q + '  -1.81%  '
$ws.Range("D7").Value = $q + '1.002'
$ws.Range("E7").Value = $q + '  -0.02%  '
$ws.Range("D8").Value = $q + '0.2557'
$ws.Range("E8").Value = $q + '  -1.87%  '
$ws.Range("D9").Value = $q + '0.06342'
$ws.Range("E9").Value = $q + '  -1.05%  '
$ws.Range("D10").Value = $q + '19.20'
$ws.Range("E10").Value = $q + '  -3.86%  '
$ws.Range("D11").Value = $q + '0.07761'
$ws.Range("E11").Value = $q + '  -0.37%  '
$ws.Range("D12").Value = $q + '4.225'
$ws.Range("E12").Value = $q + '  -2.34%  '
$ws.Range("D13").Value = $q + '1.623.05'
$ws.Range("E13").Value = $q + '  -1.92%  '
$ws.Range("D14").Value = $q + '1.843.26'
$ws.Range("E14").Value = $q + '  -1.65%  '
$ws.Range("D15").Value = $q + '0.5533'
$ws.Range("E15").Value = $q + '  +0.05%  '
$ws.Range("D16").Value = $q + '63.32'
$ws.Range("E16").Value = $q + '  -2.66%  '
$ws.Range("D17").Value = $q + '0.0' + [char]0x2085 + '7491'
$ws.Range("E17").Value = $q + '  -4.08%  '
$ws.Range("D18").Value = $q + '25.688.34'
$ws.Range("E18").Value = $q + '  -1.63%  '
$ws.Range("D19").Value = $q + '1.002'
$ws.Range("E19").Value = $q + '  -0.08%  '
$ws.Range("D20").Value = $q + '193.14'
$ws.Range("E20").Value = $q + '  -4.06%  '
$ws.Range("D21").Value = $q + '4.337'
$ws.Range("E21").Value = $q + '  -3.70%  '
$ws.Range("D22").Value = $q + '9.708'
$ws.Range("E22").Value = $q + '  -3.45%  '
$ws.Range("D23").Value = $q + '5.941'
$ws.Range("E23").Value = $q + '  -3.40%  '
$ws.Range("D24").Value = $q + '1.002'
$ws.Range("E24").Value = $q + '  -0.18%  '
$ws.Range("D25").Value = $q + '1.837'
$ws.Range("E25").Value = $q + '  -3.49%  '
$ws.Range("D26").Value = $q + '140.09'
$ws.Range("E26").Value = $q + '  -1.77%  '
$ws.Range("D27").Value = $q + '0.1253'
$ws.Range("E27").Value = $q + '  +2.48%  '
$ws.Range("D28").Value = $q + '6.694'
$ws.Range("E28").Value = $q + '  -3.52%  '
$ws.Range("D29").Value = $q + '15.39'
$ws.Range("E29").Value = $q + '  -2.52%  '
$ws.Range("D30").Value = $q + '1.230'
$ws.Range("E30").Value = $q + '  -1.28%  '
$ws.Range("D31").Value = $q + '0.04834'
$ws.Range("E31").Value = $q + '  -2.39%  '
$ws.Range("D32").Value = $q + '3.279'
$ws.Range("E32").Value = $q + '  -1.94%  '
$ws.Range("D33").Value = $q + '3.160'
$ws.Range("E33").Value = $q + '  -3.51%  '
$ws.Range("E34").Value = $q + '  -1.10%  '
$ws.Range("D35").Value = $q + '2.362'
$ws.Range("E35").Value = $q + '  -0.75%  '
$ws.Range("D36").Value = $q + '0.8876'
$ws.Range("E36").Value = $q + '  -4.27%  '
$ws.Range("D37").Value = $q + '1.120.57'
$ws.Range("E37").Value = $q + '  +0.10%  '
$ws.Range("D38").Value = $q + '2.527'
$ws.Range("E38").Value = $q + '  -2.57%  '
$ws.Range("D39").Value = $q + '0.5456'
$ws.Range("E39").Value = $q + '  -3.20%  '
$ws.Range("B40").Value = $q + 'PaxDollar'
$ws.Range("C40").Value = $q + 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = $q + '1.028'
$ws.Range("E40").Value = $q + '  +2.60%  '
$ws.Range("B41").Value = $q + 'VeChain'
$ws.Range("C41").Value = $q + 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = $q + '0.01553'
$ws.Range("E41").Value = $q + '  -1.77%  '
$ws.Range("D42").Value = $q + '5.547'
$ws.Range("E42").Value = $q + '  -1.18%  '
$ws.Range("D43").Value = $q + '0.7903'
$ws.Range("E43").Value = $q + '  -3.00%  '
$ws.Range("D44").Value = $q + '96.79'
$ws.Range("E44").Value = $q + '  -3.30%  '
$ws.Range("D45").Value = $q + '1.767.85'
$ws.Range("E45").Value = $q + '  -0.89%  '
$ws.Range("D46").Value = $q + '0.0' + [char]0x2088 + '111'
$ws.Range("E46").Value = $q + '  -9.42%  '
$ws.Range("D47").Value = $q + '0.4406'
$ws.Range("E47").Value = $q + '  -3.05%  '
$ws.Range("D48").Value = $q + '54.41'
$ws.Range("E48").Value = $q + '  -2.29%  '
$ws.Range("D49").Value = $q + '0.05092'
$ws.Range("E49").Value = $q + '  -3.53%  '
$ws.Range("D50").Value = $q + '7.510'
$ws.Range("E50").Value = $q + '  +0.23%  '
$ws.Range("D51").Value = $q + '0.9938'
$ws.Range("E51").Value = $q + '  -1.02%  '
